$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new weekly time-sheet rows (14-16) for 2026-02-09, 2026-02-10, 2026-02-11 ---
# Each new row mirrors the existing data rows: DATE | CLOCK IN (8:00 AM) | CLOCK OUT (4:00 PM) | TOTAL HOURS (8)

# 1) Copy the formatting (style) of the last existing data row (13) onto the new rows.
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D16").PasteSpecial(-4122)  # xlPasteFormats

# 2) Write the date values as plain text (not auto-converted to date serials) by
#    generating the text through a throwaway formula cell and pasting only the
#    resulting *value* into the target cell - this keeps the cell a normal text
#    (shared-string) cell instead of turning it into a date number or formula cell.
$ws.Range("F1").Formula = '="2026-02-09"'
$ws.Range("F1").Copy()
$ws.Range("A14").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("F1").Formula = '="2026-02-10"'
$ws.Range("F1").Copy()
$ws.Range("A15").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("F1").Formula = '="2026-02-11"'
$ws.Range("F1").Copy()
$ws.Range("A16").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("F1").ClearContents()

# 3) Fill in the remaining columns for the new rows.
$ws.Range("B14").Value = "8:00 AM"
$ws.Range("C14").Value = "4:00 PM"
$ws.Range("D14").Value = 8

$ws.Range("B15").Value = "8:00 AM"
$ws.Range("C15").Value = "4:00 PM"
$ws.Range("D15").Value = 8

$ws.Range("B16").Value = "8:00 AM"
$ws.Range("C16").Value = "4:00 PM"
$ws.Range("D16").Value = 8

# 4) Update the sheet's selection to cover the newly expanded table (A1:D16).
$ws.Range("A1:D16").Select() | Out-Null
